$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Rows 3 & 4: only the K/L (SalesQuantity / Turnover) figures move.
# ------------------------------------------------------------------
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 3.41

$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 3.36

# ------------------------------------------------------------------
# Rows 6 & 7: the two product rows swap places (A, E, F, G, H, I, J),
# while K/L get their own new totals. Use a scratch row (101, well
# outside the sheet's used range) plus Range.Copy so the text cells
# keep their original shared-string type and style (a plain
# .Value = "0208" assignment would coerce the text back to a number
# and lose the leading zero / change styles.xml).
# ------------------------------------------------------------------
$ws.Range("A6:A6").Copy($ws.Range("A101"))
$ws.Range("E6:E6").Copy($ws.Range("E101"))
$ws.Range("F6:F6").Copy($ws.Range("F101"))
$ws.Range("J6:J6").Copy($ws.Range("J101"))

$ws.Range("A7:A7").Copy($ws.Range("A6"))
$ws.Range("E7:E7").Copy($ws.Range("E6"))
$ws.Range("F7:F7").Copy($ws.Range("F6"))
$ws.Range("J7:J7").Copy($ws.Range("J6"))

$ws.Range("A101:A101").Copy($ws.Range("A7"))
$ws.Range("E101:E101").Copy($ws.Range("E7"))
$ws.Range("F101:F101").Copy($ws.Range("F7"))
$ws.Range("J101:J101").Copy($ws.Range("J7"))

# Drop the scratch row completely so the sheet dimension goes back to
# A1:L11 instead of staying stretched to row 101.
$ws.Range("A101").EntireRow.Delete()

# Numeric cells in rows 6 & 7 (plain numbers - safe to set directly).
$ws.Range("G6").Value = 3.7
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 50
$ws.Range("K6").Value = 8
$ws.Range("L6").Value = 13.08

$ws.Range("G7").Value = 0.85
$ws.Range("H7").Value = 0.85
$ws.Range("I7").Value = 0
$ws.Range("K7").Value = 8.414999999999999
$ws.Range("L7").Value = 6.33

# ------------------------------------------------------------------
# Rows 8-10: only K/L change.
# ------------------------------------------------------------------
$ws.Range("K8").Value = 9.295
$ws.Range("L8").Value = 13.57

$ws.Range("K9").Value = 14
$ws.Range("L9").Value = 37.28

$ws.Range("K10").Value = 45.557
$ws.Range("L10").Value = 34.89

# ------------------------------------------------------------------
# Row 11: totals.
# ------------------------------------------------------------------
$ws.Range("K11").Value = 91.267
$ws.Range("L11").Value = 115.76
